$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "wrong email" negative test-case row (row 4)
$ws.Range("A4").Value = "sdf"
$ws.Range("B4").Value = "facebook@123"
$ws.Range("C4").Value = "NEGATIVE"
$ws.Range("D4").Value = "The email address or phone number that you've entered doesn't match any account"

# Mirror the existing hyperlinked email/username cells (B2, B3) on the new row
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:facebook@123")

# Leave the selection where the author ended up after entering the new row
$ws.Range("D5").Select() | Out-Null
